# Scheduled market-price refresh for Belias_Profits: update the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for the
# rows whose crafted-item market data changed on this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1306.4028
$ws.Cells.Item(15, 9).Value = 1306.4028
$ws.Cells.Item(15, 11).Value = 3919.2084
$ws.Cells.Item(15, 13).Value = -3750.2084
$ws.Cells.Item(39, 8).Value = 449.9
$ws.Cells.Item(39, 9).Value = 133.625
$ws.Cells.Item(39, 10).Value = 1715
$ws.Cells.Item(39, 11).Value = 400.875
$ws.Cells.Item(39, 12).Value = 5145
$ws.Cells.Item(39, 13).Value = -104.875
$ws.Cells.Item(39, 14).Value = -5737
$ws.Cells.Item(68, 8).Value = 28766.334
$ws.Cells.Item(68, 10).Value = 28766.334
$ws.Cells.Item(68, 12).Value = 28766.334
$ws.Cells.Item(68, 14).Value = -30264.334
$ws.Cells.Item(71, 8).Value = 28766.334
$ws.Cells.Item(71, 10).Value = 28766.334
$ws.Cells.Item(71, 12).Value = 86299.00199999999
$ws.Cells.Item(71, 14).Value = -93787.00199999999
$ws.Cells.Item(76, 8).Value = 3100
$ws.Cells.Item(76, 9).Value = 3100
$ws.Cells.Item(76, 11).Value = 3100
$ws.Cells.Item(76, 13).Value = -2785
$ws.Cells.Item(79, 8).Value = 3100
$ws.Cells.Item(79, 9).Value = 3100
$ws.Cells.Item(79, 11).Value = 3100
$ws.Cells.Item(79, 13).Value = -2008

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 6994024
$ws.Cells.Item(45, 9).Value = 8265301
$ws.Cells.Item(45, 11).Value = 8265301
$ws.Cells.Item(45, 13).Value = -8264924
$ws.Cells.Item(74, 8).Value = 38099.418
$ws.Cells.Item(74, 9).Value = 47170.137
$ws.Cells.Item(74, 10).Value = 15926.556
$ws.Cells.Item(74, 11).Value = 47170.137
$ws.Cells.Item(74, 12).Value = 15926.556
$ws.Cells.Item(74, 13).Value = -46296.137
$ws.Cells.Item(74, 14).Value = -17674.556
$ws.Cells.Item(76, 8).Value = 21214.666
$ws.Cells.Item(76, 10).Value = 21214.666
$ws.Cells.Item(76, 12).Value = 21214.666
$ws.Cells.Item(76, 14).Value = -21890.666
$ws.Cells.Item(77, 8).Value = 38099.418
$ws.Cells.Item(77, 9).Value = 47170.137
$ws.Cells.Item(77, 10).Value = 15926.556
$ws.Cells.Item(77, 11).Value = 235850.685
$ws.Cells.Item(77, 12).Value = 79632.78
$ws.Cells.Item(77, 13).Value = -231482.685
$ws.Cells.Item(77, 14).Value = -88368.78
$ws.Cells.Item(79, 8).Value = 21214.666
$ws.Cells.Item(79, 10).Value = 21214.666
$ws.Cells.Item(79, 12).Value = 21214.666
$ws.Cells.Item(79, 14).Value = -23554.666
$ws.Cells.Item(132, 8).Value = 2686.6206
$ws.Cells.Item(132, 9).Value = 1983.4117
$ws.Cells.Item(132, 10).Value = 3682.8333
$ws.Cells.Item(132, 11).Value = 5950.2351
$ws.Cells.Item(132, 12).Value = 11048.4999
$ws.Cells.Item(132, 13).Value = -3420.2351
$ws.Cells.Item(132, 14).Value = -16108.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2064.074
$ws.Cells.Item(99, 9).Value = 1571.5385
$ws.Cells.Item(99, 10).Value = 2521.4285
$ws.Cells.Item(99, 11).Value = 1571.5385
$ws.Cells.Item(99, 12).Value = 2521.4285
$ws.Cells.Item(99, 13).Value = -73.53850000000011
$ws.Cells.Item(99, 14).Value = -5517.4285
$ws.Cells.Item(134, 8).Value = 2505764.5
$ws.Cells.Item(134, 9).Value = 3338133.2
$ws.Cells.Item(134, 10).Value = 8658.5
$ws.Cells.Item(134, 11).Value = 10014399.6
$ws.Cells.Item(134, 12).Value = 25975.5
$ws.Cells.Item(134, 13).Value = -10011864.6
$ws.Cells.Item(134, 14).Value = -31045.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1857.3
$ws.Cells.Item(31, 9).Value = 1422.3334
$ws.Cells.Item(31, 10).Value = 2213.182
$ws.Cells.Item(31, 11).Value = 1422.3334
$ws.Cells.Item(31, 12).Value = 2213.182
$ws.Cells.Item(31, 13).Value = -1127.3334
$ws.Cells.Item(31, 14).Value = -2803.182
$ws.Cells.Item(34, 8).Value = 1857.3
$ws.Cells.Item(34, 9).Value = 1422.3334
$ws.Cells.Item(34, 10).Value = 2213.182
$ws.Cells.Item(34, 11).Value = 1422.3334
$ws.Cells.Item(34, 12).Value = 2213.182
$ws.Cells.Item(34, 13).Value = -1220.3334
$ws.Cells.Item(34, 14).Value = -2617.182
$ws.Cells.Item(58, 8).Value = 7969.75
$ws.Cells.Item(58, 9).Value = 9718
$ws.Cells.Item(58, 11).Value = 9718
$ws.Cells.Item(58, 13).Value = -9515
$ws.Cells.Item(136, 8).Value = 7969.75
$ws.Cells.Item(136, 9).Value = 9718
$ws.Cells.Item(136, 11).Value = 29154
$ws.Cells.Item(136, 13).Value = -26604

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 1156384.6
$ws.Cells.Item(46, 9).Value = 900
$ws.Cells.Item(46, 10).Value = 1503030
$ws.Cells.Item(46, 11).Value = 2700
$ws.Cells.Item(46, 12).Value = 4509090
$ws.Cells.Item(46, 13).Value = -2609
$ws.Cells.Item(46, 14).Value = -4509272
$ws.Cells.Item(47, 8).Value = 37.75
$ws.Cells.Item(47, 9).Value = 37.75
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 113.25
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = 317.75
$ws.Cells.Item(47, 14).ClearContents()
$ws.Cells.Item(100, 8).Value = 3347.0588
$ws.Cells.Item(100, 10).Value = 3347.0588
$ws.Cells.Item(100, 12).Value = 10041.1764
$ws.Cells.Item(100, 14).Value = -11663.1764
$ws.Cells.Item(131, 8).Value = 940.08
$ws.Cells.Item(131, 10).Value = 940.08
$ws.Cells.Item(131, 12).Value = 2820.24
$ws.Cells.Item(131, 14).Value = -12900.24

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35, 8).Value = 4280.5
$ws.Cells.Item(35, 9).Value = 4280.5
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 4280.5
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -3944.5
$ws.Cells.Item(35, 14).ClearContents()
$ws.Cells.Item(46, 8).Value = 866.8421
$ws.Cells.Item(46, 9).Value = 892.5
$ws.Cells.Item(46, 10).Value = 848.1818
$ws.Cells.Item(46, 11).Value = 892.5
$ws.Cells.Item(46, 12).Value = 848.1818
$ws.Cells.Item(46, 13).Value = -704.5
$ws.Cells.Item(46, 14).Value = -1224.1818
$ws.Cells.Item(132, 8).Value = 5433.75
$ws.Cells.Item(132, 9).Value = 5401
$ws.Cells.Item(132, 10).Value = 5499.25
$ws.Cells.Item(132, 11).Value = 16203
$ws.Cells.Item(132, 12).Value = 16497.75
$ws.Cells.Item(132, 13).Value = -13673
$ws.Cells.Item(132, 14).Value = -21557.75
$ws.Cells.Item(136, 8).Value = 6106.2383
$ws.Cells.Item(136, 9).Value = 1306.8235
$ws.Cells.Item(136, 10).Value = 26503.75
$ws.Cells.Item(136, 11).Value = 3920.4705
$ws.Cells.Item(136, 12).Value = 79511.25
$ws.Cells.Item(136, 13).Value = -1370.4705
$ws.Cells.Item(136, 14).Value = -84611.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 7964.4
$ws.Cells.Item(32, 9).Value = 6608.6665
$ws.Cells.Item(32, 10).Value = 9998
$ws.Cells.Item(32, 11).Value = 6608.6665
$ws.Cells.Item(32, 12).Value = 9998
$ws.Cells.Item(32, 13).Value = -6291.6665
$ws.Cells.Item(32, 14).Value = -10632
$ws.Cells.Item(43, 8).Value = 9450
$ws.Cells.Item(43, 9).Value = 9500
$ws.Cells.Item(43, 10).Value = 9400
$ws.Cells.Item(43, 11).Value = 9500
$ws.Cells.Item(43, 12).Value = 9400
$ws.Cells.Item(43, 13).Value = -9351
$ws.Cells.Item(43, 14).Value = -9698

